$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.805.01'
$ws.Range("E2").Value = '  -3.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.494.31'
$ws.Range("E3").Value = '  -2.59%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.39'
$ws.Range("E5").Value = '  -2.86%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.10'
$ws.Range("E6").Value = '  -5.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.495.35'
$ws.Range("E7").Value = '  -2.56%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.489'
$ws.Range("E9").Value = '  -2.07%  '

$ws.Range("E10").Value = '  -1.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.17'
$ws.Range("E11").Value = '  -0.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.383'
$ws.Range("E12").Value = '  -2.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.079.07'
$ws.Range("E13").Value = '  -2.95%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.54'
$ws.Range("E14").Value = '  -1.76%  '

$ws.Range("E15").Value = '  +1.05%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000178'
$ws.Range("E16").Value = '  -3.53%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.517.14'
$ws.Range("E17").Value = '  -2.21%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.937.96'
$ws.Range("E18").Value = '  -3.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.13'
$ws.Range("E19").Value = '  +1.60%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.41'
$ws.Range("E20").Value = '  -1.39%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.66'
$ws.Range("E21").Value = '  -2.71%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '386.37'
$ws.Range("E22").Value = '  -2.44%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.578'
$ws.Range("E23").Value = '  -1.40%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.632.99'
$ws.Range("E24").Value = '  -2.74%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.09'
$ws.Range("E25").Value = '  -2.54%  '

$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("E27").Value = '  -3.99%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.59'
$ws.Range("E28").Value = '  -3.58%  '

$ws.Range("E29").Value = '  -0.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.41'
$ws.Range("E30").Value = '  -7.92%  '

$ws.Range("E31").Value = '  -3.45%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.21'
$ws.Range("E32").Value = '  -4.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.497.91'
$ws.Range("E33").Value = '  -2.86%  '

$ws.Range("E34").Value = '  -0.02%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.71'
$ws.Range("E35").Value = '  -3.22%  '

$ws.Range("E36").Value = '  -3.58%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.25'
$ws.Range("E37").Value = '  -2.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.56'
$ws.Range("E38").Value = '  -2.82%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.87'
$ws.Range("E39").Value = '  -1.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '167.40'
$ws.Range("E40").Value = '  -1.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0800'
$ws.Range("E41").Value = '  -4.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.07'
$ws.Range("E42").Value = '  +3.08%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.810'
$ws.Range("E43").Value = '  -3.45%  '

$ws.Range("E44").Value = '  -0.15%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.60'
$ws.Range("E45").Value = '  -3.47%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").Value = '1.20'
$ws.Range("E46").Value = '  -4.69%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").Value = '4.38'
$ws.Range("E47").Value = '  -3.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.64'
$ws.Range("E48").Value = '  -3.87%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.443.38'
$ws.Range("E49").Value = '  +0.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.86'
$ws.Range("E50").Value = '  -0.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.890'
$ws.Range("E51").Value = '  -1.83%  '
